# Append a new row (row 6) of user/vehicle data to Sheet1, extending the
# used range from A1:C5 to A1:C6 (NIK, Plat, Nama columns).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NIK is a long numeric-looking identifier; format the cell as Text first
# so Excel stores it verbatim instead of collapsing it into a double
# (same treatment the NIK values in A2:A5 already have).
$ws.Range("A6").NumberFormat = "@"
$ws.Range("A6").Value = "2345678990112444"

# Plat (vehicle plate) is left blank for this user, matching the blank
# Plat cells already present for the "Tata"/"Nia" rows.
$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value = ""

$ws.Range("C6").Value = "Tiara"
